$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be stored as text even when it looks like a number
# (e.g. "100.96"), without leaving a lingering custom number-format style on
# the cell (ClearFormats resets the style index back to the default/unset one).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = '43.620.25'
$ws.Range("E2").Value = '  +0.97%  '
$ws.Range("D3").Value = '2.245.59'
$ws.Range("E3").Value = '  +0.44%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("E5").Value = '  +0.89%  '
Set-TextValue $ws.Range("D6") '100.96'
$ws.Range("E6").Value = '  +0.17%  '
$ws.Range("E7").Value = '  -1.31%  '
Set-TextValue $ws.Range("D9") '0.554'
$ws.Range("E9").Value = '  -1.59%  '
Set-TextValue $ws.Range("D10") '36.88'
$ws.Range("E10").Value = '  -0.24%  '
Set-TextValue $ws.Range("D11") '0.0830'
$ws.Range("E11").Value = '  +0.37%  '
Set-TextValue $ws.Range("D12") '7.73'
$ws.Range("E12").Value = '  +0.56%  '
$ws.Range("E13").Value = '  -2.24%  '
$ws.Range("D14").Value = '2.585.14'
$ws.Range("E14").Value = '  +0.34%  '
Set-TextValue $ws.Range("D15") '0.854'
$ws.Range("E15").Value = '  -0.91%  '
Set-TextValue $ws.Range("D16") '14.10'
$ws.Range("E16").Value = '  -1.87%  '
$ws.Range("D17").Value = '2.241.63'
$ws.Range("E17").Value = '  +0.18%  '
$ws.Range("D18").Value = '43.495.72'
$ws.Range("E18").Value = '  +0.85%  '
Set-TextValue $ws.Range("D19") '13.52'
$ws.Range("E19").Value = '  -5.97%  '
$ws.Range("D20").Value = '0.0₃0984'
$ws.Range("E20").Value = '  +1.86%  '
Set-TextValue $ws.Range("D21") '6.56'
$ws.Range("E21").Value = '  +0.47%  '
Set-TextValue $ws.Range("D22") '65.23'
$ws.Range("E22").Value = '  -0.30%  '
$ws.Range("E23").Value = '  -0.42%  '
Set-TextValue $ws.Range("D24") '237.23'
$ws.Range("E24").Value = '  -0.38%  '
Set-TextValue $ws.Range("D25") '2.15'
$ws.Range("E25").Value = '  -0.32%  '
$ws.Range("E26").Value = '  +0.23%  '
Set-TextValue $ws.Range("D27") '10.09'
$ws.Range("E27").Value = '  +0.32%  '
$ws.Range("E28").Value = '  -2.45%  '
Set-TextValue $ws.Range("D29") '36.64'
$ws.Range("E29").Value = '  +3.30%  '
Set-TextValue $ws.Range("D30") '6.30'
$ws.Range("E30").Value = '  -1.41%  '
Set-TextValue $ws.Range("D31") '160.13'
$ws.Range("E31").Value = '  +4.38%  '
Set-TextValue $ws.Range("D32") '20.16'
$ws.Range("E32").Value = '  -1.58%  '
Set-TextValue $ws.Range("D33") '0.0853'
$ws.Range("E33").Value = '  -2.88%  '
$ws.Range("E34").Value = '  -2.46%  '
$ws.Range("E35").Value = '  -1.61%  '
$ws.Range("E36").Value = '  +7.49%  '
Set-TextValue $ws.Range("D37") '1.91'
$ws.Range("E37").Value = '  -3.25%  '
$ws.Range("E38").Value = '  -2.15%  '
Set-TextValue $ws.Range("D39") '3.75'
$ws.Range("E39").Value = '  +1.75%  '
Set-TextValue $ws.Range("D40") '4.24'
$ws.Range("E40").Value = '  -5.02%  '
Set-TextValue $ws.Range("D41") '15.52'
$ws.Range("E41").Value = '  +20.24%  '
$ws.Range("E42").Value = '  -2.30%  '
$ws.Range("E43").Value = '  +0.23%  '
$ws.Range("D44").Value = '1.796.26'
$ws.Range("E44").Value = '  -0.17%  '
Set-TextValue $ws.Range("D45") '0.199'
$ws.Range("E45").Value = '  -3.37%  '
Set-TextValue $ws.Range("D46") '82.36'
$ws.Range("E46").Value = '  -5.51%  '
Set-TextValue $ws.Range("D47") '74.73'
$ws.Range("E47").Value = '  -2.36%  '
Set-TextValue $ws.Range("D48") '5.18'
$ws.Range("E48").Value = '  -2.86%  '
Set-TextValue $ws.Range("D49") '58.58'
$ws.Range("E49").Value = '  -1.31%  '
Set-TextValue $ws.Range("D50") '103.47'
$ws.Range("E50").Value = '  +0.04%  '
Set-TextValue $ws.Range("D51") '1.67'
$ws.Range("E51").Value = '  +3.25%  '
